# "Generate Report for Handoff"
#
# b.md has progressed from "Handed back: in sync with en-US" to
# "Ready for handoff" in both locales. This refreshes the localization
# status report accordingly:
#   - Overview sheet: b.md row now shows "Ready for handoff" for both
#     zh-cn and de-de.
#   - zh-cn / de-de detail sheets: b.md row's Status becomes
#     "Ready for handoff", and a new handoff file + handoff datetime are
#     recorded for the Latest Handoff File / Latest Handoff Datetime
#     columns (the hyperlink text is updated to match).

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"

# --- Overview sheet: row 3 is b.md ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus

# --- zh-cn sheet: row 3 is b.md ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B3").Value = $newStatus
$zhcn.Range("C3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("D3").Value = "2016-03-09 10:32:39"
foreach ($hl in $zhcn.Hyperlinks) {
    if ($hl.Range.Address() -eq '$C$3') {
        $hl.TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
    }
}

# --- de-de sheet: row 3 is b.md ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B3").Value = $newStatus
$dede.Range("C3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("D3").Value = "2016-03-09 10:32:43"
foreach ($hl in $dede.Hyperlinks) {
    if ($hl.Range.Address() -eq '$C$3') {
        $hl.TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
    }
}
